$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment"); this shifts C:Q left to B:P,
# matching the diff's column-shift for all data rows and the new
# dimension A1:P33.
$ws.Columns.Item(2).Delete()

# Append ".jamais.jamais" to every header cell in row 1 except "Country"
# (now in columns B through P after the shift).
for ($col = 2; $col -le 16; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Text + ".jamais.jamais"
}
